# Finished Most of Receive
# Implement Frame for Actions Taken and show Actions Taken in Report.
#
# Populates the (previously empty) "Reports" sheet with a header row and
# one data row, and narrows column C to fit the new "Report" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column C (Report) from 88.8 to ~68.4 characters wide.
$ws.Columns.Item(3).ColumnWidth = 67.5

# Header row
$ws.Range("A1").Value = "Request ID"
$ws.Range("B1").Value = "Report Timestamp"
$ws.Range("C1").Value = "Report"

# First report row
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "2025-04-25 19:30:06"
$ws.Range("C2").Value = "John Smith found battery 2. Now John Smith is Confident"
